# Added support for merge import
#
# This updates the "famhist_*" columns (J:Q) of the single data sheet so the
# columns are reordered/relabeled, the Yes/No answers in rows 2-6 are
# corrected to reflect the new column meanings, and one new label ("fdsa")
# is introduced. It also moves the "text" number-format style that used to
# sit on J2 over to L2, and updates the active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- selection -----------------------------------------------------------
$ws.Range("Q15").Select()

# --- style bookkeeping -----------------------------------------------------
# In the original file J2 carried the custom "text" number format (s="1").
# In the edited file that same format instead belongs to L2, and J2 goes
# back to the default (unstyled) cell.
$ws.Range("J2").ClearFormats()
$ws.Range("L2").NumberFormat = "@"

# --- header row (row 1): famhist_* column labels --------------------------
$ws.Range("J1").Value = "famhist_none"
$ws.Range("K1").Value = "famhist_deaf"
$ws.Range("L1").Value = "famhist_cardiomyopathy"
$ws.Range("M1").Value = "famhist_encephalopathy"
$ws.Range("N1").Value = "famhist_diabmell"
$ws.Range("O1").Value = "famhist_cardiovasc"
$ws.Range("Q1").Value = "famhist_unknown"
# P1 (famhist_malignancy) is unchanged.

# --- row 2 ------------------------------------------------------------
$ws.Range("J2").Value = "No"
$ws.Range("K2").Value = "No"
$ws.Range("L2").Value = "Yes"
$ws.Range("M2").Value = "Yes"
$ws.Range("N2").Value = "Yes"
$ws.Range("O2").Value = "No"
$ws.Range("P2").Value = "No"
$ws.Range("Q2").Value = "No"

# --- row 3 ------------------------------------------------------------
$ws.Range("J3").Value = "fdsa"
$ws.Range("K3").Value = "Yes"
$ws.Range("L3").Value = "Yes"
$ws.Range("M3").Value = 3
$ws.Range("N3").Value = "No"
$ws.Range("O3").Value = "No"
$ws.Range("P3").Value = "No"
$ws.Range("Q3").Value = "No"

# --- row 4 ------------------------------------------------------------
$ws.Range("J4").Value = "Yes"
$ws.Range("N4").Value = "No"
$ws.Range("O4").Value = "No"
$ws.Range("P4").Value = "No"
$ws.Range("Q4").Value = "No"

# --- row 5 ------------------------------------------------------------
$ws.Range("J5").Value = 4
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = "Yes"
$ws.Range("Q5").Value = "No"

# --- row 6 ------------------------------------------------------------
$ws.Range("N6").Value = "No"
$ws.Range("O6").Value = "No"
$ws.Range("P6").Value = "No"
$ws.Range("Q6").Value = "Yes"
